$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously empty cells with new values (Sprint 4 final metrics)
$ws.Range("E3").Value = 5
$ws.Range("E5").Value = 3
$ws.Range("E7").Value = "1.  Government publish opportunity.  2.  Contractor view published opportunities              3.  Government edit opportunity"

# Adjust row 7 height to fit the new longer text
$ws.Rows.Item(7).RowHeight = 121.5

# Update the view - scroll so row 6 is at the top, and select E8
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("E8").Select()
